# Update MIAPPE templates to the new ontology version.
#
# Sheet "isa_template": the TAGS block (rows 13-15) drops the MIAPPE-specific
# "Observation Unit" tag and adds a new "plant material" tag (FOODON-based)
# as an additional column.
#
# Sheet "observation_unit_sample": the example data row's Plant Anatomical
# Entity / Plant structure development stage accession-number URLs are
# updated to the new OLS4-style PURL-resolver links.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("isa_template")

# Row 13 ("Tags" header row): shift the tag names left, dropping
# "Observation Unit" and appending the new "plant material" tag.
$ws1.Range("C13").Value = "study"
$ws1.Range("D13").Value = "growth"
$ws1.Range("E13").Value = "plant material"

# Row 14 ("Tags Term Accession Number" row): shift the accession numbers
# left to match, dropping the MIAPPE_0069 accession and appending the new
# FOODON accession for "plant material".
$ws1.Range("C14").Value = "http://purl.obolibrary.org/obo/NCIT_C63536"
$ws1.Range("D14").Value = "http://purl.obolibrary.org/obo/GO_0040007"
$ws1.Range("E14").Value = "http://purl.obolibrary.org/obo/FOODON_00004331"

# Row 15 ("Tags Term Source REF" row): add the new tag's term source.
$ws1.Range("E15").Value = "FOODON"

$ws2 = $wb.Worksheets.Item("observation_unit_sample")

# Update the example row's ontology term accession URLs to the new
# OLS4-resolver format.
$ws2.Range("M2").Value = "https://www.ebi.ac.uk/ols4/ontologies/po/classes/http%253A%252F%252Fpurl.obolibrary.org%252Fobo%252FPO_0025094"
$ws2.Range("P2").Value = "https://www.ebi.ac.uk/ols4/ontologies/po/classes/http%253A%252F%252Fpurl.obolibrary.org%252Fobo%252FPO_0000003"
